$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$newValue = "mixed or unspecified population"

# Column B addresses (contiguous row blocks) whose value changes from either
# "regular, several popualtions, or unspecified" or "high risk" to
# "mixed or unspecified population"
$ranges = @(
    "B3:B34",
    "B37:B58",
    "B64",
    "B66:B73",
    "B75:B96",
    "B98:B117",
    "B127:B133",
    "B135:B148",
    "B150:B168",
    "B172:B214",
    "B235:B239",
    "B245:B272",
    "B276:B308",
    "B310:B317",
    "B338:B351",
    "B357:B358",
    "B366:B406",
    "B408:B445",
    "B447:B453",
    "B457:B475",
    "B477:B485",
    "B487:B495",
    "B500:B507",
    "B513:B515",
    "B526:B528",
    "B531:B537",
    "B540:B589",
    "B596:B602",
    "B607:B650"
)

foreach ($addr in $ranges) {
    $ws.Range($addr).Value = $newValue
}
